$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 28937.75
$ws.Range("I98").Value = 689
$ws.Range("K98").Value = 689
$ws.Range("M98").Value = 809
$ws.Range("H100").Value = 1933.2667
$ws.Range("I100").Value = 1345.3636
$ws.Range("J100").Value = 3550
$ws.Range("K100").Value = 1345.3636
$ws.Range("L100").Value = 3550
$ws.Range("M100").Value = -804.3635999999999
$ws.Range("N100").Value = -4632
$ws.Range("H122").Value = 28937.75
$ws.Range("I122").Value = 689
$ws.Range("K122").Value = 2067
$ws.Range("M122").Value = 383
$ws.Range("H132").Value = 39568.44
$ws.Range("I132").Value = 5838.4707
$ws.Range("J132").Value = 111244.625
$ws.Range("K132").Value = 17515.4121
$ws.Range("L132").Value = 333733.875
$ws.Range("M132").Value = -14985.4121
$ws.Range("N132").Value = -338793.875
$ws.Range("H137").Value = 2082974.5
$ws.Range("I137").Value = 6993912
$ws.Range("K137").Value = 20981736
$ws.Range("M137").Value = -20979186
$ws.Range("H138").Value = 2624.2222
$ws.Range("I138").Value = 2295.4285
$ws.Range("J138").Value = 2833.4546
$ws.Range("K138").Value = 6886.2855
$ws.Range("L138").Value = 8500.363799999999
$ws.Range("M138").Value = -1746.2855
$ws.Range("N138").Value = -18780.3638
$ws.Range("H141").Value = 3050.5264
$ws.Range("I141").Value = 1637.0834
$ws.Range("K141").Value = 4911.2502
$ws.Range("M141").Value = 268.7497999999996

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2536.818
$ws.Range("I2").Value = 2536.818
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2536.818
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2423.818
$ws.Range("N2").ClearContents()
$ws.Range("H4").Value = 285.33334
$ws.Range("I4").Value = 253.33333
$ws.Range("J4").Value = 317.33334
$ws.Range("K4").Value = 253.33333
$ws.Range("L4").Value = 317.33334
$ws.Range("M4").Value = -137.33333
$ws.Range("N4").Value = -549.33334
$ws.Range("H32").Value = 3849.72
$ws.Range("I32").Value = 1724.2667
$ws.Range("K32").Value = 1724.2667
$ws.Range("M32").Value = -1437.2667
$ws.Range("H61").Value = 2764.077
$ws.Range("I61").Value = 1491.1428
$ws.Range("J61").Value = 4249.1665
$ws.Range("K61").Value = 1491.1428
$ws.Range("L61").Value = 4249.1665
$ws.Range("M61").Value = -1279.1428
$ws.Range("N61").Value = -4673.1665
$ws.Range("H74").Value = 1695.8914
$ws.Range("I74").Value = 1335.4865
$ws.Range("K74").Value = 1335.4865
$ws.Range("M74").Value = -461.4865
$ws.Range("H77").Value = 1695.8914
$ws.Range("I77").Value = 1335.4865
$ws.Range("K77").Value = 6677.4325
$ws.Range("M77").Value = -2309.4325
$ws.Range("H110").Value = 1937.3636
$ws.Range("I110").Value = 1851.375
$ws.Range("J110").Value = 2166.6667
$ws.Range("K110").Value = 1851.375
$ws.Range("L110").Value = 2166.6667
$ws.Range("M110").Value = 193.625
$ws.Range("N110").Value = -6256.6667
$ws.Range("H116").Value = 2536.818
$ws.Range("I116").Value = 2536.818
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2536.818
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -242.8180000000002
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 33336612
$ws.Range("I132").Value = 71430310
$ws.Range("J132").Value = 4624.25
$ws.Range("K132").Value = 214290930
$ws.Range("L132").Value = 13872.75
$ws.Range("M132").Value = -214288400
$ws.Range("N132").Value = -18932.75
$ws.Range("H136").Value = 2764.077
$ws.Range("I136").Value = 1491.1428
$ws.Range("J136").Value = 4249.1665
$ws.Range("K136").Value = 4473.428400000001
$ws.Range("L136").Value = 12747.4995
$ws.Range("M136").Value = -1923.428400000001
$ws.Range("N136").Value = -17847.4995

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2536.818
$ws.Range("I3").Value = 2536.818
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2536.818
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2422.818
$ws.Range("N3").ClearContents()
$ws.Range("H86").Value = 2225.9
$ws.Range("I86").Value = 2119.875
$ws.Range("K86").Value = 2119.875
$ws.Range("M86").Value = -996.875
$ws.Range("H89").Value = 2225.9
$ws.Range("I89").Value = 2119.875
$ws.Range("K89").Value = 10599.375
$ws.Range("M89").Value = -4983.375
$ws.Range("H105").Value = 3652.7334
$ws.Range("I105").Value = 2742
$ws.Range("J105").Value = 4108.1
$ws.Range("K105").Value = 2742
$ws.Range("L105").Value = 4108.1
$ws.Range("M105").Value = -995
$ws.Range("N105").Value = -7602.1

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 89536.75
$ws.Range("I132").Value = 1169.3
$ws.Range("J132").Value = 236815.83
$ws.Range("K132").Value = 3507.9
$ws.Range("L132").Value = 710447.49
$ws.Range("M132").Value = -977.8999999999996
$ws.Range("N132").Value = -715507.49
$ws.Range("H134").Value = 2583250.2
$ws.Range("I134").Value = 4696834.5
$ws.Range("J134").Value = 469665.66
$ws.Range("K134").Value = 14090503.5
$ws.Range("L134").Value = 1408996.98
$ws.Range("M134").Value = -14087968.5
$ws.Range("N134").Value = -1414066.98

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1000004
$ws.Range("I14").Value = 1000004
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1000004
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -999836
$ws.Range("N14").ClearContents()
$ws.Range("H102").Value = 906
$ws.Range("I102").Value = 812
$ws.Range("K102").Value = 812
$ws.Range("M102").Value = 810
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 40005316
$ws.Range("I132").Value = 90915650
$ws.Range("J132").Value = 4341.5
$ws.Range("K132").Value = 272746950
$ws.Range("L132").Value = 13024.5
$ws.Range("M132").Value = -272744420
$ws.Range("N132").Value = -18084.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3572.5833
$ws.Range("I7").Value = 2630.8572
$ws.Range("J7").Value = 4891
$ws.Range("K7").Value = 2630.8572
$ws.Range("L7").Value = 4891
$ws.Range("M7").Value = -2518.8572
$ws.Range("N7").Value = -5115
$ws.Range("H126").Value = 3572.5833
$ws.Range("I126").Value = 2630.8572
$ws.Range("J126").Value = 4891
$ws.Range("K126").Value = 7892.571599999999
$ws.Range("L126").Value = 14673
$ws.Range("M126").Value = -5422.571599999999
$ws.Range("N126").Value = -19613
$ws.Range("H128").Value = 28429
$ws.Range("J128").Value = 28429
$ws.Range("L128").Value = 28429
$ws.Range("N128").Value = -38389
